# Updates Price (D) and Volume(1h) (E) columns for the cryptos table on Sheet1,
# reflecting the latest pull of market data (see commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "39.752.38"
$ws.Cells.Item(2, 5).Value = "  +1.43%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.213.56"
$ws.Cells.Item(3, 5).Value = "  +0.78%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.07%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "292.43"
$ws.Cells.Item(5, 5).Value = "  -0.56%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "86.77"
$ws.Cells.Item(6, 5).Value = "  +7.61%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  +1.35%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.05%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.472"
$ws.Cells.Item(9, 5).Value = "  +1.45%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "30.28"
$ws.Cells.Item(10, 5).Value = "  +4.61%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0786"
$ws.Cells.Item(11, 5).Value = "  +2.62%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "47.38"
$ws.Cells.Item(12, 5).Value = "  +1.44%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.109"
$ws.Cells.Item(13, 5).Value = "  +1.90%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "6.34"
$ws.Cells.Item(14, 5).Value = "  +2.47%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.553.34"
$ws.Cells.Item(15, 5).Value = "  +0.81%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.04"
$ws.Cells.Item(16, 5).Value = "  +1.36%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.197.16"
$ws.Cells.Item(17, 5).Value = "  +0.10%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.728"
$ws.Cells.Item(18, 5).Value = "  +3.00%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "39.695.48"
$ws.Cells.Item(19, 5).Value = "  +1.58%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "11.47"
$ws.Cells.Item(20, 5).Value = "  +12.33%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "0.0₃0880"
$ws.Cells.Item(21, 5).Value = "  +1.55%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.80"
$ws.Cells.Item(22, 5).Value = "  +2.19%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "65.74"
$ws.Cells.Item(23, 5).Value = "  +1.99%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "235.68"
$ws.Cells.Item(24, 5).Value = "  +4.65%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  +0.03%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.47"
$ws.Cells.Item(26, 5).Value = "  +3.37%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "1.84"
$ws.Cells.Item(27, 5).Value = "  +2.68%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "22.75"
$ws.Cells.Item(28, 5).Value = "  +1.54%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +1.60%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "9.27"
$ws.Cells.Item(30, 5).Value = "  +3.10%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "32.81"
$ws.Cells.Item(31, 5).Value = "  +4.37%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "152.36"
$ws.Cells.Item(32, 5).Value = "  +2.21%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.998"
$ws.Cells.Item(33, 5).Value = "  -0.18%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "4.93"
$ws.Cells.Item(34, 5).Value = "  +3.24%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.0718"
$ws.Cells.Item(35, 5).Value = "  +3.86%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.38"
$ws.Cells.Item(36, 5).Value = "  +1.68%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +2.51%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +7.25%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "16.00"
$ws.Cells.Item(39, 5).Value = "  +5.21%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.0987"
$ws.Cells.Item(40, 5).Value = "  +3.67%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.70"
$ws.Cells.Item(41, 5).Value = "  +3.62%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "2.079.55"
$ws.Cells.Item(42, 5).Value = "  +9.61%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.79"
$ws.Cells.Item(43, 5).Value = "  +5.69%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +5.67%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.0269"
$ws.Cells.Item(45, 5).Value = "  +4.17%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.01"
$ws.Cells.Item(46, 5).Value = "  +11.65%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "17.78"
$ws.Cells.Item(47, 5).Value = "  +12.40%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "2.61"
$ws.Cells.Item(48, 5).Value = "  +0.89%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "2.423.17"
$ws.Cells.Item(49, 5).Value = "  +0.78%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "70.86"
$ws.Cells.Item(50, 5).Value = "  -0.26%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "89.12"
$ws.Cells.Item(51, 5).Value = "  +2.52%  "
